$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H121").Value = 600
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H137").Value = 1178.8214
$ws.Range("I137").Value = 1019.625
$ws.Range("J137").Value = 1391.0834
$ws.Range("K137").Value = 3058.875
$ws.Range("L137").Value = 4173.2502
$ws.Range("M137").Value = -508.875
$ws.Range("N137").Value = -9273.2502
$ws.Range("H141").Value = 2752.2122
$ws.Range("I141").Value = 1391.9048
$ws.Range("J141").Value = 5132.75
$ws.Range("K141").Value = 4175.7144
$ws.Range("L141").Value = 15398.25
$ws.Range("M141").Value = 1004.2856
$ws.Range("N141").Value = -25758.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1808.91
$ws.Range("I32").Value = 1332.8353
$ws.Range("J32").Value = 4506.6665
$ws.Range("K32").Value = 1332.8353
$ws.Range("L32").Value = 4506.6665
$ws.Range("M32").Value = -1045.8353
$ws.Range("N32").Value = -5080.6665
$ws.Range("H45").Value = 1347.9445
$ws.Range("I45").Value = 1255.4445
$ws.Range("J45").Value = 1440.4445
$ws.Range("K45").Value = 1255.4445
$ws.Range("L45").Value = 1440.4445
$ws.Range("M45").Value = -878.4445000000001
$ws.Range("N45").Value = -2194.4445
$ws.Range("H61").Value = 1429.9354
$ws.Range("I61").Value = 805.3333
$ws.Range("J61").Value = 3571.4285
$ws.Range("K61").Value = 805.3333
$ws.Range("L61").Value = 3571.4285
$ws.Range("M61").Value = -593.3333
$ws.Range("N61").Value = -3995.4285
$ws.Range("H88").Value = 1118714.1
$ws.Range("I88").Value = 2508001.5
$ws.Range("J88").Value = 7284.2
$ws.Range("K88").Value = 2508001.5
$ws.Range("L88").Value = 7284.2
$ws.Range("M88").Value = -2507595.5
$ws.Range("N88").Value = -8096.2
$ws.Range("H91").Value = 1118714.1
$ws.Range("I91").Value = 2508001.5
$ws.Range("J91").Value = 7284.2
$ws.Range("K91").Value = 2508001.5
$ws.Range("L91").Value = 7284.2
$ws.Range("M91").Value = -2506597.5
$ws.Range("N91").Value = -10092.2
$ws.Range("H132").Value = 4923.121
$ws.Range("I132").Value = 6676.5264
$ws.Range("J132").Value = 2543.5
$ws.Range("K132").Value = 20029.5792
$ws.Range("L132").Value = 7630.5
$ws.Range("M132").Value = -17499.5792
$ws.Range("N132").Value = -12690.5
$ws.Range("H136").Value = 1429.9354
$ws.Range("I136").Value = 805.3333
$ws.Range("J136").Value = 3571.4285
$ws.Range("K136").Value = 2415.9999
$ws.Range("L136").Value = 10714.2855
$ws.Range("M136").Value = 134.0001000000002
$ws.Range("N136").Value = -15814.2855

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 631.5714
$ws.Range("I107").Value = 661.8333
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 661.8333
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1258.1667
$ws.Range("N107").Value = -4290

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 37039376
$ws.Range("I62").Value = 2005
$ws.Range("J62").Value = 83336090
$ws.Range("K62").Value = 2005
$ws.Range("L62").Value = 83336090
$ws.Range("M62").Value = -1381
$ws.Range("N62").Value = -83337338
$ws.Range("H65").Value = 37039376
$ws.Range("I65").Value = 2005
$ws.Range("J65").Value = 83336090
$ws.Range("K65").Value = 10025
$ws.Range("L65").Value = 416680450
$ws.Range("M65").Value = -6905
$ws.Range("N65").Value = -416686690
$ws.Range("H86").Value = 62502220
$ws.Range("I86").Value = 142858060
$ws.Range("K86").Value = 142858060
$ws.Range("M86").Value = -142856937
$ws.Range("H89").Value = 62502220
$ws.Range("I89").Value = 142858060
$ws.Range("K89").Value = 714290300
$ws.Range("M89").Value = -714284684

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 24025.176
$ws.Range("I120").Value = 8010
$ws.Range("J120").Value = 27457
$ws.Range("K120").Value = 24030
$ws.Range("L120").Value = 82371
$ws.Range("M120").Value = -19192
$ws.Range("N120").Value = -92047
$ws.Range("H131").Value = 2472346.2
$ws.Range("J131").Value = 3268877.2
$ws.Range("L131").Value = 9806631.600000001
$ws.Range("N131").Value = -9816711.600000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8333.333000000001
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 8333.333000000001
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -59984
$ws.Range("H132").Value = 50654.094
$ws.Range("I132").Value = 79340.30499999999
$ws.Range("J132").Value = 4039
$ws.Range("K132").Value = 238020.915
$ws.Range("L132").Value = 12117
$ws.Range("M132").Value = -235490.915
$ws.Range("N132").Value = -17177

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I64").Value = 30000
$ws.Range("K64").Value = 30000
$ws.Range("M64").Value = -29775
$ws.Range("I67").Value = 30000
$ws.Range("K67").Value = 30000
$ws.Range("M67").Value = -29220
$ws.Range("H136").Value = 6604.125
$ws.Range("I136").Value = 8166.6
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 24499.8
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -21949.8
$ws.Range("N136").Value = -17100

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 19800
$ws.Range("J64").Value = 19800
$ws.Range("L64").Value = 19800
$ws.Range("N64").Value = -20296
$ws.Range("H67").Value = 19800
$ws.Range("J67").Value = 19800
$ws.Range("L67").Value = 19800
$ws.Range("N67").Value = -21516
$ws.Range("H81").Value = 2160
$ws.Range("I81").Value = 1650
$ws.Range("J81").Value = 2670
$ws.Range("K81").Value = 3300
$ws.Range("L81").Value = 5340
$ws.Range("M81").Value = -2239
$ws.Range("N81").Value = -7462
$ws.Range("H84").Value = 2160
$ws.Range("I84").Value = 1650
$ws.Range("J84").Value = 2670
$ws.Range("K84").Value = 16500
$ws.Range("L84").Value = 26700
$ws.Range("M84").Value = -11196
$ws.Range("N84").Value = -37308
